# Update column G ("K") values in the save_data sheet.
# The commit regenerates the "K" column (previously tracked as "Strike#")
# from the underlying simulation logs. The new values below are the
# regenerated K counts for each data row (rows 2-67).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 0
    3 = 0
    4 = 1
    5 = 1
    6 = 0
    7 = 0
    8 = 1
    9 = 0
    10 = 0
    11 = 2
    12 = 2
    13 = 1
    14 = 0
    15 = 1
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 2
    21 = 0
    22 = 3
    23 = 1
    24 = 0
    25 = 0
    26 = 0
    27 = 1
    28 = 1
    29 = 1
    30 = 1
    31 = 3
    32 = 1
    33 = 0
    34 = 0
    35 = 1
    36 = 2
    37 = 1
    38 = 0
    39 = 2
    40 = 1
    41 = 0
    42 = 0
    43 = 2
    44 = 0
    45 = 0
    46 = 0
    47 = 2
    48 = 2
    49 = 0
    50 = 1
    51 = 1
    52 = 1
    53 = 0
    54 = 1
    55 = 0
    56 = 3
    57 = 2
    58 = 0
    59 = 1
    60 = 2
    61 = 2
    62 = 1
    63 = 0
    64 = 3
    65 = 1
    66 = 1
    67 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
